$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.848.54'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.634.87'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.17'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.32'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.664.24'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.23'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.336'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.099.63'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.855.55'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.85'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.661.82'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '348.26'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.52'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.33'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.79'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.418'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.161'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0801'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.09'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.28'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +7.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.97'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.57'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.77'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.29%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.99'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.13'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.62'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.836'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.70'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.29%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '278.14'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.77%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0985'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.55'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0529'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0229'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.983.12'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.65'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.50%  '
